$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (bold header style, same style as other headers like H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF)
$iValues = @(5, 5, 5, 8, 8, 6, 6, 8, 9, 5, 9, 8, 9)
$jValues = @(5, 6, 6, 9, 8, 6, 7, 8, 9, 6, 9, 8, 9)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
